$d = $word.ActiveDocument

# Update the date line in the first paragraph
$d.Content.Find.Execute("2025-06-24 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-25 Wednesday", 2) | Out-Null

# Update table cells individually (handles duplicate source strings mapping to distinct targets)
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "71÷3=23, 2"

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "48÷9=5, 3"

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "90÷2=45, 0"

$cell = $t.Cell(1, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "87÷4=21, 3"

$cell = $t.Cell(1, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "35÷9=3, 8"

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "37÷6=6, 1"

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "23÷2=11, 1"

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "26÷3=8, 2"

$cell = $t.Cell(5, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "14÷8=1, 6"

$cell = $t.Cell(5, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "17÷5=3, 2"

$cell = $t.Cell(9, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "44÷2=22, 0"

$cell = $t.Cell(9, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "20÷8=2, 4"

$cell = $t.Cell(9, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "37÷9=4, 1"

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "71÷6=11, 5"

$cell = $t.Cell(9, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "44÷2=22, 0"

$cell = $t.Cell(13, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "92÷4=23, 0"

$cell = $t.Cell(13, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "82÷3=27, 1"

$cell = $t.Cell(13, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "58÷9=6, 4"

$cell = $t.Cell(13, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "68÷8=8, 4"

$cell = $t.Cell(13, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "24÷2=12, 0"

$cell = $t.Cell(17, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "59÷8=7, 3"

$cell = $t.Cell(17, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "41÷7=5, 6"

$cell = $t.Cell(17, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "96÷7=13, 5"

$cell = $t.Cell(17, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "95÷2=47, 1"

$cell = $t.Cell(17, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "54÷5=10, 4"
